$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking
# strings like "318.42" are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "48.199.93"
$ws.Range("E2").Value = "  +0.03%  "

$ws.Range("D3").Value = "2.501.25"
$ws.Range("E3").Value = "  -1.20%  "

$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "318.42"
$ws.Range("E5").Value = "  -1.84%  "

$ws.Range("D6").Value = "106.16"
$ws.Range("E6").Value = "  -2.56%  "

$ws.Range("D7").Value = "0.519"
$ws.Range("E7").Value = "  -1.59%  "

$ws.Range("E8").Value = "  -0.11%  "

$ws.Range("E9").Value = "  -3.51%  "

$ws.Range("D10").Value = "38.93"
$ws.Range("E10").Value = "  -4.24%  "

$ws.Range("D11").Value = "20.29"
$ws.Range("E11").Value = "  -1.09%  "

$ws.Range("D12").Value = "0.0804"
$ws.Range("E12").Value = "  -2.93%  "

$ws.Range("E13").Value = "  +0.11%  "

$ws.Range("D14").Value = "7.12"
$ws.Range("E14").Value = "  -2.59%  "

$ws.Range("D15").Value = "2.892.52"
$ws.Range("E15").Value = "  -1.25%  "

$ws.Range("D16").Value = "2.501.85"
$ws.Range("E16").Value = "  -1.41%  "

$ws.Range("D17").Value = "0.831"
$ws.Range("E17").Value = "  -3.48%  "

$ws.Range("D18").Value = "48.056.35"
$ws.Range("E18").Value = "  +0.02%  "

$ws.Range("D19").Value = "3.00"
$ws.Range("E19").Value = "  +11.51%  "

$ws.Range("D20").Value = "12.85"
$ws.Range("E20").Value = "  -3.38%  "

$ws.Range("D21").Value = "6.59"
$ws.Range("E21").Value = "  -1.01%  "

$ws.Range("D22").Value = "0.0₃0932"
$ws.Range("E22").Value = "  -2.10%  "

$ws.Range("D23").Value = "71.17"
$ws.Range("E23").Value = "  -1.61%  "

$ws.Range("D24").Value = "268.39"
$ws.Range("E24").Value = "  -0.53%  "

$ws.Range("D25").Value = "2.52"
$ws.Range("E25").Value = "  -2.75%  "

$ws.Range("E26").Value = "  +0.23%  "

$ws.Range("D27").Value = "25.81"
$ws.Range("E27").Value = "  -1.68%  "

$ws.Range("D28").Value = "2.20"
$ws.Range("E28").Value = "  -0.43%  "

$ws.Range("E29").Value = "  -4.03%  "

$ws.Range("D30").Value = "0.140"
$ws.Range("E30").Value = "  -3.22%  "

$ws.Range("D31").Value = "34.71"
$ws.Range("E31").Value = "  -2.82%  "

$ws.Range("D32").Value = "49.36"
$ws.Range("E32").Value = "  -0.84%  "

$ws.Range("B33").Value = "Celestia"
$ws.Range("C33").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D33").Value = "19.17"
$ws.Range("E33").Value = "  -3.62%  "

$ws.Range("B34").Value = "FirstDigitalUSD"
$ws.Range("C34").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.21%  "

$ws.Range("D35").Value = "5.30"
$ws.Range("E35").Value = "  -2.40%  "

$ws.Range("E36").Value = "  -2.77%  "

$ws.Range("D37").Value = "1.95"
$ws.Range("E37").Value = "  -2.41%  "

$ws.Range("D38").Value = "4.61"
$ws.Range("E38").Value = "  -3.54%  "

$ws.Range("E39").Value = "  -4.44%  "

$ws.Range("D40").Value = "123.30"
$ws.Range("E40").Value = "  +3.12%  "

$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").Value = "0.110"
$ws.Range("E41").Value = "  -1.87%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "22.25"
$ws.Range("E42").Value = "  -0.80%  "

$ws.Range("E43").Value = "  +1.66%  "

$ws.Range("D44").Value = "0.0303"
$ws.Range("E44").Value = "  +0.70%  "

$ws.Range("D45").Value = "2.003.53"
$ws.Range("E45").Value = "  -0.52%  "

$ws.Range("D46").Value = "3.16"
$ws.Range("E46").Value = "  +0.64%  "

$ws.Range("E47").Value = "  +1.99%  "

$ws.Range("E48").Value = "  -2.65%  "

$ws.Range("E49").Value = "  -2.46%  "

$ws.Range("D50").Value = "5.20"
$ws.Range("E50").Value = "  -1.19%  "

$ws.Range("D51").Value = "78.93"
$ws.Range("E51").Value = "  -1.07%  "
